# Applies the "Feature/52266 +: [Review intern] improve review intern" edit:
#  - appends 4 new status rows (19-22) to the ReviewDetails sample data
#  - adds the corresponding new Status strings (PmReviewed / HrApproved / ReOpen / Rejected)
#  - widens a handful of columns to fit the new/longer data
#  - nudges the active selection the way the author's last save left it

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. New data rows (19-22), mirroring the existing row layout (A = Id-2 index
#    column with the bordered/centered style already used by A2:A18, B = Id,
#    etc.). Column letters map to the header row as follows:
#      A index(=Id-2)  B Id  C CreationTime  D CreatorUserId
#      E LastModificationTime  F LastModifierUserId  G IsDeleted
#      J ReviewId  K InternshipId  L ReviewerId  M CurrentLevel  N NewLevel
#      O Status  P Note  Q Type  S Salary  U RateStar
# ---------------------------------------------------------------------------

$creationTime = "2022-12-27 23:28:48.2374196"
$modTime = "2022-12-27 23:29:52.7488269"
$note = "HR đánh giá chung:`n    ok`n"

$rows = @(
    @{ Row=19; A=17; B=18; K=6;  M="Intern_1"; N="Intern_2"; O="PmReviewed" },
    @{ Row=20; A=18; B=19; K=7;  M="Intern_3"; N="Intern_3"; O="HrApproved" },
    @{ Row=21; A=19; B=20; K=12; M="Intern_1"; N="Intern_2"; O="ReOpen" },
    @{ Row=22; A=20; B=21; K=12; M="Intern_1"; N="Intern_2"; O="Rejected" }
)

foreach ($r in $rows) {
    $row = $r.Row

    $ws.Cells.Item($row, 1).Value = $r.A          # A - index
    $ws.Cells.Item(17, 1).Copy() | Out-Null        # reuse the existing bordered/centered style (A2:A18)
    $ws.Cells.Item($row, 1).PasteSpecial(-4122) | Out-Null   # xlPasteFormats

    $ws.Cells.Item($row, 2).Value = $r.B          # B - Id
    $ws.Cells.Item($row, 3).Value = $creationTime # C - CreationTime
    $ws.Cells.Item($row, 4).Value = 1             # D - CreatorUserId
    $ws.Cells.Item($row, 5).Value = $modTime      # E - LastModificationTime
    $ws.Cells.Item($row, 6).Value = 1             # F - LastModifierUserId
    $ws.Cells.Item($row, 7).Value = $false        # G - IsDeleted

    $ws.Cells.Item($row, 10).Value = 3            # J - ReviewId
    $ws.Cells.Item($row, 11).Value = $r.K         # K - InternshipId
    $ws.Cells.Item($row, 12).Value = 3            # L - ReviewerId
    $ws.Cells.Item($row, 13).Value = $r.M         # M - CurrentLevel
    $ws.Cells.Item($row, 14).Value = $r.N         # N - NewLevel
    $ws.Cells.Item($row, 15).Value = $r.O         # O - Status
    $ws.Cells.Item($row, 16).Value = $note        # P - Note
    $ws.Cells.Item($row, 17).Value = "Internship" # Q - Type

    $ws.Cells.Item($row, 19).Value = 2000000      # S - Salary
    $ws.Cells.Item($row, 21).Value = 4            # U - RateStar
}

# ---------------------------------------------------------------------------
# 2. Column widths - widened to fit the longer values now present
#    (ColumnWidth is expressed in "characters"; the saved OOXML width ends up
#    offset by the standard 5px/MDW padding, so we back the values off by
#    5/7 to land on the same rendered width the author's Excel produced).
# ---------------------------------------------------------------------------
$ws.Columns.Item(3).ColumnWidth = 25.428571428571427
$ws.Columns.Item(5).ColumnWidth = 25.428571428571427
$ws.Columns.Item(9).ColumnWidth = 14.0
$ws.Columns.Item(10).ColumnWidth = 9.714285714285714
$ws.Columns.Item(11).ColumnWidth = 12.571428571428571
$ws.Columns.Item(12).ColumnWidth = 11.857142857142858
$ws.Columns.Item(13).ColumnWidth = 11.714285714285714
$ws.Columns.Item(14).ColumnWidth = 11.142857142857142
$ws.Columns.Item(15).ColumnWidth = 9.142857142857142
$ws.Columns.Item(16).ColumnWidth = 25.285714285714285
$ws.Columns.Item(17).ColumnWidth = 10.857142857142858

# ---------------------------------------------------------------------------
# 3. Selection left where the author's last save left it.
# ---------------------------------------------------------------------------
$ws.Range("I24").Select() | Out-Null
